$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-81)
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22)
for ($row = 2; $row -le 81; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
